# Generate Report for Archive
#
# 1) Update status text from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview sheet columns E/F, and the
#    "Status" column (C) on the zh-cn / de-de per-language sheets).
# 2) Narrow the "Status" columns (Overview!E:F, zh-cn!C, de-de!C) to
#    their new, narrower width now that the text is shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$overviewRows = $overview.UsedRange.Rows.Count
for ($r = 2; $r -le $overviewRows; $r++) {
    if ($overview.Cells.Item($r, 5).Value() -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value() -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value = $newStatus
    }
}

# --- zh-cn / de-de sheets: column C holds the status ---
foreach ($ws in @($zhcn, $dede)) {
    $usedRows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $usedRows; $r++) {
        if ($ws.Cells.Item($r, 3).Value() -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
}

# --- Narrow the status columns to match the shorter text ---
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
